# trafo_id -> gridnode_id refactor
# The "buildings" sheet has a header row (row 1) with one column per field.
# Column W holds the header "trafo_id" - rename it to "gridnode_id".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("W1").Value = "gridnode_id"

# Reflect the on-screen selection/scroll state captured at save time.
$ws.Range("X6").Select()

Write-Output "renamed W1 header to:"
Write-Output $ws.Range("W1").Value()
